# MIN_2009 worksheet: add season record columns (Wins/Losses/Ties)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting used by the existing header row (bold, centered,
# bordered) from the last header cell (AC1) onto the three new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Header labels for the new columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (87-76-0) for every player row
$wins = 87
$losses = 76
$ties = 0

for ($row = 2; $row -le 44; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # AD
    $ws.Cells.Item($row, 31).Value = $losses  # AE
    $ws.Cells.Item($row, 32).Value = $ties    # AF
}
